$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 800
$ws.Range("I12").Value = 800
$ws.Range("K12").Value = 800
$ws.Range("M12").Value = -630

$ws.Range("H40").Value = 1818.5652
$ws.Range("I40").Value = 1513.1428
$ws.Range("J40").Value = 2293.6667
$ws.Range("K40").Value = 1513.1428
$ws.Range("L40").Value = 2293.6667
$ws.Range("M40").Value = -1338.1428
$ws.Range("N40").Value = -2643.6667

$ws.Range("H64").Value = 5497.5
$ws.Range("I64").Value = 5495
$ws.Range("J64").Value = 5500
$ws.Range("K64").Value = 5495
$ws.Range("L64").Value = 5500
$ws.Range("M64").Value = -5247
$ws.Range("N64").Value = -5996

$ws.Range("H67").Value = 5497.5
$ws.Range("I67").Value = 5495
$ws.Range("J67").Value = 5500
$ws.Range("K67").Value = 5495
$ws.Range("L67").Value = 5500
$ws.Range("M67").Value = -4637
$ws.Range("N67").Value = -7216

$ws.Range("H96").Value = 52000
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -14746

$ws.Range("H98").Value = 927.8570999999999
$ws.Range("I98").Value = 927.8570999999999
$ws.Range("K98").Value = 927.8570999999999
$ws.Range("M98").Value = 570.1429000000001

$ws.Range("H106").Value = 1470
$ws.Range("I106").Value = 1470
$ws.Range("K106").Value = 1470
$ws.Range("M106").Value = -839

$ws.Range("H122").Value = 927.8570999999999
$ws.Range("I122").Value = 927.8570999999999
$ws.Range("K122").Value = 2783.5713
$ws.Range("M122").Value = -333.5712999999996

$ws.Range("H131").Value = 11740
$ws.Range("I131").Value = 11740
$ws.Range("K131").Value = 35220
$ws.Range("M131").Value = -30180

$ws.Range("H132").Value = 2916.1667
$ws.Range("I132").Value = 2916.1667
$ws.Range("K132").Value = 8748.500100000001
$ws.Range("M132").Value = -6218.500100000001

$ws.Range("H134").Value = 35000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140

$ws.Range("H135").Value = 365.46155
$ws.Range("I135").Value = 365.46155
$ws.Range("K135").Value = 3289.15395
$ws.Range("M135").Value = -754.1539499999999

$ws.Range("H141").Value = 2430.7646
$ws.Range("I141").Value = 2195.1428
$ws.Range("K141").Value = 6585.428400000001
$ws.Range("M141").Value = -1405.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2752844.2
$ws.Range("I32").Value = 3687979.8
$ws.Range("J32").Value = 778668.75
$ws.Range("K32").Value = 3687979.8
$ws.Range("L32").Value = 778668.75
$ws.Range("M32").Value = -3687692.8
$ws.Range("N32").Value = -779242.75

$ws.Range("H110").Value = 4112385.8
$ws.Range("I110").Value = 5286967.5
$ws.Range("J110").Value = 1349.5
$ws.Range("K110").Value = 5286967.5
$ws.Range("L110").Value = 1349.5
$ws.Range("M110").Value = -5284922.5
$ws.Range("N110").Value = -5439.5

$ws.Range("H122").Value = 1584.9
$ws.Range("I122").Value = 1276.4117
$ws.Range("K122").Value = 3829.2351
$ws.Range("M122").Value = -1379.2351

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 31013.285
$ws.Range("I56").Value = 7093
$ws.Range("J56").Value = 35000
$ws.Range("K56").Value = 7093
$ws.Range("L56").Value = 35000
$ws.Range("M56").Value = -6248
$ws.Range("N56").Value = -36690

$ws.Range("H58").Value = 1932.1765
$ws.Range("I58").Value = 1775.1875
$ws.Range("K58").Value = 1775.1875
$ws.Range("M58").Value = -1572.1875

$ws.Range("H62").Value = 2744.6667
$ws.Range("J62").Value = 2744.6667
$ws.Range("L62").Value = 2744.6667
$ws.Range("N62").Value = -3992.6667

$ws.Range("H65").Value = 2744.6667
$ws.Range("J65").Value = 2744.6667
$ws.Range("L65").Value = 13723.3335
$ws.Range("N65").Value = -19963.3335

$ws.Range("H136").Value = 1932.1765
$ws.Range("I136").Value = 1775.1875
$ws.Range("K136").Value = 5325.5625
$ws.Range("M136").Value = -2775.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12223069
$ws.Range("I4").Value = 12223069
$ws.Range("K4").Value = 36669207
$ws.Range("M4").Value = -36669095

$ws.Range("H33").Value = 660
$ws.Range("I33").Value = 270.5
$ws.Range("J33").Value = 1049.5
$ws.Range("K33").Value = 1623
$ws.Range("L33").Value = 6297
$ws.Range("M33").Value = -1340
$ws.Range("N33").Value = -6863

$ws.Range("H69").Value = 2500
$ws.Range("J69").Value = 2500
$ws.Range("L69").Value = 7500
$ws.Range("N69").Value = -9122

$ws.Range("H72").Value = 2500
$ws.Range("J72").Value = 2500
$ws.Range("L72").Value = 22500
$ws.Range("N72").Value = -30612

$ws.Range("H98").Value = 4000
$ws.Range("J98").Value = 4000
$ws.Range("L98").Value = 12000
$ws.Range("N98").Value = -14996

$ws.Range("H103").Value = 1049.1428
$ws.Range("I103").Value = 56
$ws.Range("J103").Value = 1446.4
$ws.Range("K103").Value = 168
$ws.Range("L103").Value = 4339.200000000001
$ws.Range("M103").Value = 711
$ws.Range("N103").Value = -6097.200000000001

$ws.Range("H137").Value = 1765.3334
$ws.Range("I137").Value = 1318.6
$ws.Range("J137").Value = 3999
$ws.Range("K137").Value = 3955.8
$ws.Range("L137").Value = 11997
$ws.Range("M137").Value = 1144.2
$ws.Range("N137").Value = -22197

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 220.3
$ws.Range("I31").Value = 220.3
$ws.Range("K31").Value = 220.3
$ws.Range("M31").Value = 71.69999999999999

$ws.Range("H33").Value = 30040000
$ws.Range("I33").Value = 80000
$ws.Range("J33").Value = 60000000
$ws.Range("K33").Value = 80000
$ws.Range("L33").Value = 60000000
$ws.Range("M33").Value = -79748
$ws.Range("N33").Value = -60000504

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H37").Value = 220.3
$ws.Range("I37").Value = 220.3
$ws.Range("K37").Value = 220.3
$ws.Range("M37").Value = 56.69999999999999

$ws.Range("H52").Value = 25030
$ws.Range("I52").Value = 25030
$ws.Range("K52").Value = 25030
$ws.Range("M52").Value = -24771

$ws.Range("H54").Value = 11536.25
$ws.Range("I54").Value = 6000
$ws.Range("J54").Value = 12327.143
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 12327.143
$ws.Range("M54").Value = -5610
$ws.Range("N54").Value = -13107.143

$ws.Range("H70").Value = 5066.375
$ws.Range("I70").Value = 4902.8
$ws.Range("J70").Value = 5339
$ws.Range("K70").Value = 4902.8
$ws.Range("L70").Value = 5339
$ws.Range("M70").Value = -4632.8
$ws.Range("N70").Value = -5879

$ws.Range("H73").Value = 5066.375
$ws.Range("I73").Value = 4902.8
$ws.Range("J73").Value = 5339
$ws.Range("K73").Value = 4902.8
$ws.Range("L73").Value = 5339
$ws.Range("M73").Value = -3966.8
$ws.Range("N73").Value = -7211

$ws.Range("H80").Value = 2808.875
$ws.Range("I80").Value = 2759
$ws.Range("J80").Value = 2838.8
$ws.Range("K80").Value = 2759
$ws.Range("L80").Value = 2838.8
$ws.Range("M80").Value = -1761
$ws.Range("N80").Value = -4834.8

$ws.Range("H83").Value = 2808.875
$ws.Range("I83").Value = 2759
$ws.Range("J83").Value = 2838.8
$ws.Range("K83").Value = 13795
$ws.Range("L83").Value = 14194
$ws.Range("M83").Value = -8803
$ws.Range("N83").Value = -24178

$ws.Range("H122").Value = 1594.4
$ws.Range("I122").Value = 1594.4
$ws.Range("K122").Value = 4783.200000000001
$ws.Range("M122").Value = -2333.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22000.666
$ws.Range("J7").Value = 6000
$ws.Range("L7").Value = 6000
$ws.Range("N7").Value = -6224

$ws.Range("H126").Value = 22000.666
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19690.572
$ws.Range("I41").Value = 19632.334
$ws.Range("K41").Value = 19632.334
$ws.Range("M41").Value = -19242.334

$ws.Range("H81").Value = 1002547.3
$ws.Range("I81").Value = 1163.3334
$ws.Range("J81").Value = 1431711.9
$ws.Range("K81").Value = 2326.6668
$ws.Range("L81").Value = 2863423.8
$ws.Range("M81").Value = -1265.6668
$ws.Range("N81").Value = -2865545.8

$ws.Range("H84").Value = 1002547.3
$ws.Range("I84").Value = 1163.3334
$ws.Range("J84").Value = 1431711.9
$ws.Range("K84").Value = 11633.334
$ws.Range("L84").Value = 14317119
$ws.Range("M84").Value = -6329.333999999999
$ws.Range("N84").Value = -14327727

$ws.Range("H126").Value = 4520.9
$ws.Range("I126").Value = 4465.6924
$ws.Range("K126").Value = 13397.0772
$ws.Range("M126").Value = -10927.0772
